# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.880.32"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.626.83"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.51"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.64"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.25"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.851.28"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.627.08"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.62"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.863.67"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.72"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.23"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.54"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.44"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0498"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.902"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.129.56"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.547"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.47"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.762.05"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.25"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("E48").Value = "  +4.08%  "
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"
$ws.Range("E51").Value = "  +1.89%  "
